{"js": "// The author split \"ze? Nev\u00edte. Tak j\u00e1 v\u00e1m to \u0159eknu.\" so that:\n//   - \"Nev\u00edte.\" became \"Nev\u00edte? Lidi!\"\n//   - the phrase \" Tak j\u00e1 v\u00e1m to \u0159eknu.\" moved to AFTER the existing\n//     \"_GoBack\" bookmark (it used to precede it).\n// Net visible text change:\n//   \"...pen\u00edze? Nev\u00edte. Tak j\u00e1 v\u00e1m to \u0159eknu. A j\u00e1 tu dnes...\"\n//   -> \"...pen\u00edze? Nev\u00edte? Lidi! Tak j\u00e1 v\u00e1m to \u0159eknu. A j\u00e1 tu dnes...\"\n// with the bookmark staying anchored right after \"Lidi!\".\n\nconst body = context.document.body;\n\n// Step 1: turn \"ze? Nev\u00edte. Tak j\u00e1 v\u00e1m to \u0159eknu.\" into \"ze? Nev\u00edte? Lidi!\"\n// (this is the text that currently sits right before the bookmark).\nconst beforeBookmark = body.search(\"ze? Nev\u00edte. Tak j\u00e1 v\u00e1m to \u0159eknu.\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nbeforeBookmark.load(\"items\");\nawait context.sync();\n\nif (beforeBookmark.items.length === 0) {\n  throw new Error(\"Could not find the target sentence to update.\");\n}\nbeforeBookmark.items[0].insertText(\"ze? Nev\u00edte? Lidi!\", \"Replace\");\nawait context.sync();\n\n// Step 2: re-insert \" Tak j\u00e1 v\u00e1m to \u0159eknu.\" right after the bookmark,\n// i.e. immediately before \"A j\u00e1 tu dnes jsem od toho\" (note: no leading\n// space in the search text so the match starts exactly at \"A\", leaving\n// the run's pre-existing leading space untouched ahead of our insertion).\nconst afterBookmark = body.search(\"A j\u00e1 tu dnes jsem od toho\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nafterBookmark.load(\"items\");\nawait context.sync();\n\nif (afterBookmark.items.length === 0) {\n  throw new Error(\"Could not find the anchor text after the bookmark.\");\n}\nafterBookmark.items[0].insertText(\"Tak j\u00e1 v\u00e1m to \u0159eknu. \", \"Start\");\nawait context.sync();\n", "ps1": "# The author split \"ze? Nev\u00edte. Tak j\u00e1 v\u00e1m to \u0159eknu.\" so that:\n#   - \"Nev\u00edte.\" became \"Nev\u00edte? Lidi!\"\n#   - the phrase \" Tak j\u00e1 v\u00e1m to \u0159eknu.\" moved to AFTER the existing\n#     \"_GoBack\" bookmark (it used to precede it).\n# Net visible text change:\n#   \"...pen\u00edze? Nev\u00edte. Tak j\u00e1 v\u00e1m to \u0159eknu. A j\u00e1 tu dnes...\"\n#   -> \"...pen\u00edze? Nev\u00edte? Lidi! Tak j\u00e1 v\u00e1m to \u0159eknu. A j\u00e1 tu dnes...\"\n# with the bookmark staying anchored right after \"Lidi!\".\n\n$d = $word.ActiveDocument\n\n# Step 1: turn \"ze? Nev\u00edte. Tak j\u00e1 v\u00e1m to \u0159eknu.\" into \"ze? Nev\u00edte? Lidi!\"\n# (this is the text that currently sits right before the bookmark).\n$r1 = $d.Content\n$found1 = $r1.Find.Execute(\n    \"ze? Nev\u00edte. Tak j\u00e1 v\u00e1m to \u0159eknu.\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"ze? Nev\u00edte? Lidi!\", 2\n)\nif (-not $found1) {\n    throw \"Could not find the target sentence to update.\"\n}\n\n# Step 2: re-insert \" Tak j\u00e1 v\u00e1m to \u0159eknu.\" right after the bookmark, i.e.\n# immediately before \"A j\u00e1 tu dnes jsem od toho\" (note: no leading space in\n# the search text so the match starts exactly at \"A\", leaving the run's\n# pre-existing leading space untouched ahead of our insertion). Collapsing\n# the found range to its start before inserting keeps the bookmark's\n# position (which sits right before this text) intact.\n$r2 = $d.Content\n$found2 = $r2.Find.Execute(\"A j\u00e1 tu dnes jsem od toho\")\nif (-not $found2) {\n    throw \"Could not find the anchor text after the bookmark.\"\n}\n$r2.Collapse(1)\n$r2.InsertBefore(\"Tak j\u00e1 v\u00e1m to \u0159eknu. \")\n"}
